# Work_breakdown_structure.xlsx edit
# - Drop the trailing "." from the two task-title cells (B3, B4 on the
#   single worksheet "Folha1").
# - Move the active cell / selection from C8 to B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "1.1 Search for Code Smells"
$ws.Range("B4").Value = "1.2 Search for Gof Design Patterns"

$ws.Range("B5").Select()
